$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 3 corresponds to file "metrics_sim_with_priors.json" (B3).
# Update the metric values to the corrected figures per the commit:
# "Correcting Relevance Markers Walker (2018) - Wolters (2018)"

$ws.Range("C3").Value = 0.7244094488188977
$ws.Range("D3").Value = 0.9068241469816273
$ws.Range("H3").Value = 0.593343669250646
$ws.Range("I3").Value = 0.07898692504026751
$ws.Range("J3").Value = 0.6246719160104987
$ws.Range("K3").Value = 4142.304461942257

$ws.Range("Q3").Value = 43
$ws.Range("R3").Value = 188
$ws.Range("S3").Value = 1347
$ws.Range("T3").Value = 4881
$ws.Range("U3").Value = 10614
$ws.Range("V3").Value = 47570
$ws.Range("W3").Value = 47425
$ws.Range("X3").Value = 46266
$ws.Range("Y3").Value = 42732
$ws.Range("Z3").Value = 36999

$ws.Range("AF3").Value = 0.999097
$ws.Range("AG3").Value = 0.996051
$ws.Range("AH3").Value = 0.971709
$ws.Range("AI3").Value = 0.897486
$ws.Range("AJ3").Value = 0.777078
